$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: best_params updates
$ws.Range("F2").Value = "{'max_depth': 50, 'n_estimators': 50}"
$ws.Range("K2").Value = "{'activation': 'leaky_relu', 'b_random_vec_range': [0, 10], 'lam': 2, 'n_layer': 64, 'n_nodes': 128, 'random_seed': 856, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Row 3: rmse updates
$ws.Range("F3").Value = 0.07633399699235413
$ws.Range("G3").Value = 0.07172586620849958
$ws.Range("H3").Value = 0.1124061764308056
$ws.Range("K3").Value = 0.04526424254504183

# Row 4: r2 updates
$ws.Range("F4").Value = 0.8004907225253259
$ws.Range("G4").Value = 0.8253679590716827
$ws.Range("H4").Value = 0.5613696538801605
$ws.Range("K4").Value = 0.9302770881995229

# Row 5: mape updates
$ws.Range("F5").Value = 19.23281259306575
$ws.Range("G5").Value = 22.18864652085109
$ws.Range("H5").Value = 84.72932397939235
$ws.Range("K5").Value = 13.51650447913631

$wb.Save()
